# edit.ps1 - apply the "more fixes to the slides" commit to IntroductionToSoftwareClass.pptx
#
# Summary of changes (see diff):
#   1. presentation.xml gains an (empty) p15:sldGuideLst extension -- this is a
#      PowerPoint view-state artifact with no content implications; the object
#      model does not expose a way to create it (Presentation.Guides is an
#      inert stub in this host), so it is attempted defensively but otherwise
#      skipped.
#   2. Slide 10 (C example): "    /* error check here */" becomes
#      "    " + "/* " + "add error checking " + "here */" (4 runs).
#   3. Slide 8 (R example):
#        a. normAutofit fontScale/lnSpcReduction reset to a bare <a:normAutofit/>
#        b. the stray "https://github.com/jasoncoposky/training/tree/master/2014-oss/day-09"
#           URL text that leaked into the "( 10, 110 ... )" call is removed.
#   4. Slide 9 (Python example): "    # error check " becomes
#      "    # " + "add error checking " (2 runs), the trailing "here" run is untouched.

$p = $ppt.ActivePresentation
$nbsp = [char]0x00A0

# ---------------------------------------------------------------------------
# 1. presentation.xml -- p:extLst / p15:sldGuideLst
#    Best-effort only: touch Presentation.Guides so that if the host ever
#    wires this up, the empty guide list extension gets persisted on save.
# ---------------------------------------------------------------------------
try {
    $null = $p.Guides
} catch {
}

# ---------------------------------------------------------------------------
# 2. Slide 10 -- "Examples: A Function in C"
#    Paragraph 2 of the content placeholder: "    /* error check here */"
#    -> "    " / "/* " / "add error checking " / "here */"
# ---------------------------------------------------------------------------
$slideC = $p.Slides.Item(10)
$shapeC = $slideC.Shapes.Item(2)
$trC = $shapeC.TextFrame.TextRange

# "error check " (positions 53-64 of the full text) -> "add error checking "
$runC1 = $trC.Characters(53, 12)
$runC1.Text = "add error checking "

# Split the leading "    " (NBSP SP NBSP SP) away from "/* "
$runC2 = $trC.Characters(46, 4)
$runC2.Text = "$nbsp $nbsp "

# Split "/* " away from what is now "add error checking here */"
$runC3 = $trC.Characters(50, 3)
$runC3.Text = "/* "

# ---------------------------------------------------------------------------
# 3. Slide 8 -- "Examples: A Function in R"
# ---------------------------------------------------------------------------
$slideR = $p.Slides.Item(8)
$shapeR = $slideR.Shapes.Item(2)

# 3a. Drop the autofit shrink amounts -> <a:normAutofit/>
$shapeR.TextFrame.AutoSize = 2   # ppAutoSizeTextToFitShape

# 3b. Fix "( 10, 110https://github.com/jasoncoposky/training/tree/master/2014-oss/day-09 )"
#     -> "( 10, 110 )"
$trR = $shapeR.TextFrame.TextRange
$spanR = $trR.Characters(173, 79)
$spanR.Text = "( 10, 110 )"
# Split "( " away from "10, 110 )"
$runR1 = $trR.Characters(173, 2)
$runR1.Text = "( "

# ---------------------------------------------------------------------------
# 4. Slide 9 -- "Examples: A Function in Python"
#    Paragraph 2: "    # error check " -> "    # " / "add error checking "
#    (the following "here" run is left untouched)
# ---------------------------------------------------------------------------
$slidePy = $p.Slides.Item(9)
$shapePy = $slidePy.Shapes.Item(2)
$trPy = $shapePy.TextFrame.TextRange

# "error check " (positions 37-48) -> "add error checking "
$runPy1 = $trPy.Characters(37, 12)
$runPy1.Text = "add error checking "

# Split the leading "    # " (NBSP SP NBSP SP # SP) away from the new text
$runPy2 = $trPy.Characters(31, 6)
$runPy2.Text = "$nbsp $nbsp # "
